$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-like text cells need to be pre-formatted as Text so Excel
# does not auto-convert "NN%" strings into numeric percentage values.
$pctCells = @("H3", "H8", "H16", "H19", "H21", "H22", "H23", "H24", "H32", "H34", "H39", "H40")
foreach ($addr in $pctCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-27 18:18:43"
$ws.Range("E3").Value = "2026-02-27 18:18:46"
$ws.Range("H3").Value = "37%"
$ws.Range("O3").Value = "4.8 °C"
$ws.Range("E4").Value = "2026-02-27 18:18:49"
$ws.Range("J4").Value = "1024.6 hPa"
$ws.Range("E5").Value = "2026-02-27 18:18:51"
$ws.Range("O5").Value = "5.2 °C"
$ws.Range("E6").Value = "2026-02-27 18:18:54"
$ws.Range("J6").Value = "1024.5 hPa"
$ws.Range("O6").Value = "11.2 °C"
$ws.Range("E7").Value = "2026-02-27 18:18:57"
$ws.Range("O7").Value = "11.5 °C"
$ws.Range("E8").Value = "2026-02-27 18:19:00"
$ws.Range("H8").Value = "60%"
$ws.Range("N8").Value = "9.0 °C 17:57 TU"
$ws.Range("O8").Value = "12.3 °C"
$ws.Range("E9").Value = "2026-02-27 18:19:02"
$ws.Range("E10").Value = "2026-02-27 18:19:05"
$ws.Range("E11").Value = "2026-02-27 18:19:08"
$ws.Range("O11").Value = "8.8 °C"
$ws.Range("E12").Value = "2026-02-27 18:19:10"
$ws.Range("E13").Value = "2026-02-27 18:19:12"
$ws.Range("E14").Value = "2026-02-27 18:19:13"
$ws.Range("E15").Value = "2026-02-27 18:19:16"
$ws.Range("E16").Value = "2026-02-27 18:19:18"
$ws.Range("H16").Value = "41%"
$ws.Range("E17").Value = "2026-02-27 18:19:21"
$ws.Range("K17").Value = "16.6 MJ/m2"
$ws.Range("N17").Value = "5.2 °C 17:58 TU"
$ws.Range("O17").Value = "7.9 °C"
$ws.Range("E18").Value = "2026-02-27 18:19:23"
$ws.Range("J18").Value = "1024.7 hPa"
$ws.Range("E19").Value = "2026-02-27 18:19:25"
$ws.Range("H19").Value = "63%"
$ws.Range("E20").Value = "2026-02-27 18:19:27"
$ws.Range("E21").Value = "2026-02-27 18:19:30"
$ws.Range("H21").Value = "60%"
$ws.Range("J21").Value = "1024.6 hPa"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-27 18:19:33"
$ws.Range("H22").Value = "50%"
$ws.Range("E23").Value = "2026-02-27 18:19:35"
$ws.Range("H23").Value = "39%"
$ws.Range("E24").Value = "2026-02-27 18:19:38"
$ws.Range("H24").Value = "75%"
$ws.Range("J24").Value = "1023.8 hPa"
$ws.Range("O24").Value = "10.4 °C"
$ws.Range("E25").Value = "2026-02-27 18:19:41"
$ws.Range("O25").Value = "6.4 °C"
$ws.Range("E26").Value = "2026-02-27 18:19:44"
$ws.Range("O26").Value = "10.7 °C"
$ws.Range("E27").Value = "2026-02-27 18:19:46"
$ws.Range("O27").Value = "5.9 °C"
$ws.Range("E28").Value = "2026-02-27 18:19:49"
$ws.Range("E29").Value = "2026-02-27 18:19:52"
$ws.Range("O29").Value = "11.7 °C"
$ws.Range("E30").Value = "2026-02-27 18:19:55"
$ws.Range("J30").Value = "1024.6 hPa"
$ws.Range("E31").Value = "2026-02-27 18:19:57"
$ws.Range("J31").Value = "1024.3 hPa"
$ws.Range("E32").Value = "2026-02-27 18:20:00"
$ws.Range("H32").Value = "60%"
$ws.Range("O32").Value = "8.1 °C"
$ws.Range("E33").Value = "2026-02-27 18:20:03"
$ws.Range("J33").Value = "1024.0 hPa"
$ws.Range("O33").Value = "8.7 °C"
$ws.Range("E34").Value = "2026-02-27 18:20:06"
$ws.Range("H34").Value = "45%"
$ws.Range("O34").Value = "4.9 °C"
$ws.Range("E35").Value = "2026-02-27 18:20:08"
$ws.Range("J35").Value = "1022.7 hPa"
$ws.Range("O35").Value = "12.3 °C"
$ws.Range("E36").Value = "2026-02-27 18:20:11"
$ws.Range("E37").Value = "2026-02-27 18:20:13"
$ws.Range("J37").Value = "1025.1 hPa"
$ws.Range("O37").Value = "8.2 °C"
$ws.Range("E38").Value = "2026-02-27 18:20:16"
$ws.Range("O38").Value = "10.4 °C"
$ws.Range("E39").Value = "2026-02-27 18:20:19"
$ws.Range("H39").Value = "32%"
$ws.Range("E40").Value = "2026-02-27 18:20:21"
$ws.Range("H40").Value = "67%"
$ws.Range("J40").Value = "1025.1 hPa"
$ws.Range("E41").Value = "2026-02-27 18:20:24"
$ws.Range("J41").Value = "1024.8 hPa"
$ws.Range("E42").Value = "2026-02-27 18:20:27"
$ws.Range("O42").Value = "11.6 °C"
$ws.Range("E43").Value = "2026-02-27 18:20:29"
$ws.Range("O43").Value = "9.4 °C"
$ws.Range("E44").Value = "2026-02-27 18:20:32"
$ws.Range("E45").Value = "2026-02-27 18:20:34"
$ws.Range("E46").Value = "2026-02-27 18:20:37"
$ws.Range("O46").Value = "10.8 °C"
